$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data rows per latest scrape
# Row 2
$dStyle = $ws.Range("D2").Style
$ws.Range("D2").Value = "'36.241.03"
$ws.Range("D2").Style = $dStyle
$ws.Range("E2").Value = '  -3.05%  '

# Row 3
$dStyle = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.961.16"
$ws.Range("D3").Style = $dStyle
$ws.Range("E3").Value = '  -4.11%  '

# Row 4
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$dStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'243.41"
$ws.Range("D5").Style = $dStyle
$ws.Range("E5").Value = '  -3.45%  '

# Row 6
$dStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = "'0.619"
$ws.Range("D6").Style = $dStyle
$ws.Range("E6").Value = '  -4.91%  '

# Row 7
$dStyle = $ws.Range("D7").Style
$ws.Range("D7").Value = "'57.20"
$ws.Range("D7").Style = $dStyle
$ws.Range("E7").Value = '  -13.36%  '

# Row 8
$ws.Range("E8").Value = '  +0.13%  '

# Row 9
$dStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.367"
$ws.Range("D9").Style = $dStyle
$ws.Range("E9").Value = '  -8.85%  '

# Row 10
$dStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = "'55.73"
$ws.Range("D10").Style = $dStyle
$ws.Range("E10").Value = '  -6.30%  '

# Row 11
$dStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.0849"
$ws.Range("D11").Style = $dStyle
$ws.Range("E11").Value = '  +5.19%  '

# Row 12
$ws.Range("E12").Value = '  -0.48%  '

# Row 13
$dStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.833"
$ws.Range("D13").Style = $dStyle
$ws.Range("E13").Value = '  -8.92%  '

# Row 14
$dStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'21.77"
$ws.Range("D14").Style = $dStyle
$ws.Range("E14").Value = '  -7.50%  '

# Row 15
$dStyle = $ws.Range("D15").Style
$ws.Range("D15").Value = "'2.248.80"
$ws.Range("D15").Style = $dStyle
$ws.Range("E15").Value = '  -4.11%  '

# Row 16
$dStyle = $ws.Range("D16").Style
$ws.Range("D16").Value = "'13.48"
$ws.Range("D16").Style = $dStyle
$ws.Range("E16").Value = '  -9.00%  '

# Row 17
$dStyle = $ws.Range("D17").Style
$ws.Range("D17").Value = "'5.32"
$ws.Range("D17").Style = $dStyle
$ws.Range("E17").Value = '  -6.87%  '

# Row 18
$dStyle = $ws.Range("D18").Style
$ws.Range("D18").Value = "'1.969.59"
$ws.Range("D18").Style = $dStyle
$ws.Range("E18").Value = '  -3.78%  '

# Row 19
$dStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'36.086.57"
$ws.Range("D19").Style = $dStyle
$ws.Range("E19").Value = '  -3.15%  '

# Row 20
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$dStyle = $ws.Range("D20").Style
$ws.Range("D20").Value = "'70.75"
$ws.Range("D20").Style = $dStyle
$ws.Range("E20").Value = '  -3.35%  '

# Row 21
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$dStyle = $ws.Range("D21").Style
$ws.Range("D21").Value = "'0.0₃0881"
$ws.Range("D21").Style = $dStyle
$ws.Range("E21").Value = '  -1.84%  '

# Row 22
$dStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'231.21"
$ws.Range("D22").Style = $dStyle
$ws.Range("E22").Value = '  -3.28%  '

# Row 23
$dStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'5.12"
$ws.Range("D23").Style = $dStyle
$ws.Range("E23").Value = '  -7.12%  '

# Row 24
$ws.Range("E24").Value = '  +0.09%  '

# Row 25
$dStyle = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.48"
$ws.Range("D25").Style = $dStyle
$ws.Range("E25").Value = '  -4.52%  '

# Row 26
$dStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'2.27"
$ws.Range("D26").Style = $dStyle
$ws.Range("E26").Value = '  -4.46%  '

# Row 27
$dStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = "'9.58"
$ws.Range("D27").Style = $dStyle
$ws.Range("E27").Value = '  -4.95%  '

# Row 28
$dStyle = $ws.Range("D28").Style
$ws.Range("D28").Value = "'164.45"
$ws.Range("D28").Style = $dStyle
$ws.Range("E28").Value = '  +1.44%  '

# Row 29
$dStyle = $ws.Range("D29").Style
$ws.Range("D29").Value = "'19.80"
$ws.Range("D29").Style = $dStyle
$ws.Range("E29").Value = '  -1.65%  '

# Row 30
$dStyle = $ws.Range("D30").Style
$ws.Range("D30").Value = "'0.128"
$ws.Range("D30").Style = $dStyle
$ws.Range("E30").Value = '  -2.28%  '

# Row 31
$dStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = "'0.118"
$ws.Range("D31").Style = $dStyle
$ws.Range("E31").Value = '  -3.65%  '

# Row 32
$dStyle = $ws.Range("D32").Style
$ws.Range("D32").Value = "'1.15"
$ws.Range("D32").Style = $dStyle
$ws.Range("E32").Value = '  -3.65%  '

# Row 33
$dStyle = $ws.Range("D33").Style
$ws.Range("D33").Value = "'4.74"
$ws.Range("D33").Style = $dStyle
$ws.Range("E33").Value = '  -7.96%  '

# Row 34
$dStyle = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.0639"
$ws.Range("D34").Style = $dStyle
$ws.Range("E34").Value = '  +1.31%  '

# Row 35
$dStyle = $ws.Range("D35").Style
$ws.Range("D35").Value = "'4.34"
$ws.Range("D35").Style = $dStyle
$ws.Range("E35").Value = '  -6.96%  '

# Row 36
$ws.Range("E36").Value = '  +0.15%  '

# Row 37
$ws.Range("E37").Value = '  -2.00%  '

# Row 38
$dStyle = $ws.Range("D38").Style
$ws.Range("D38").Value = "'5.94"
$ws.Range("D38").Style = $dStyle
$ws.Range("E38").Value = '  -7.32%  '

# Row 39
$dStyle = $ws.Range("D39").Style
$ws.Range("D39").Value = "'2.14"
$ws.Range("D39").Style = $dStyle
$ws.Range("E39").Value = '  -10.14%  '

# Row 40
$ws.Range("E40").Value = '  -3.46%  '

# Row 41
$dStyle = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.0956"
$ws.Range("D41").Style = $dStyle
$ws.Range("E41").Value = '  -5.97%  '

# Row 42
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$dStyle = $ws.Range("D42").Style
$ws.Range("D42").Value = "'2.88"
$ws.Range("D42").Style = $dStyle
$ws.Range("E42").Value = '  -5.23%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$dStyle = $ws.Range("D43").Style
$ws.Range("D43").Value = "'1.18"
$ws.Range("D43").Style = $dStyle
$ws.Range("E43").Value = '  -9.30%  '

# Row 44
$dStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.0210"
$ws.Range("D44").Style = $dStyle
$ws.Range("E44").Value = '  -4.27%  '

# Row 45
$dStyle = $ws.Range("D45").Style
$ws.Range("D45").Value = "'1.06"
$ws.Range("D45").Style = $dStyle
$ws.Range("E45").Value = '  -9.33%  '

# Row 46
$dStyle = $ws.Range("D46").Style
$ws.Range("D46").Value = "'15.69"
$ws.Range("D46").Style = $dStyle
$ws.Range("E46").Value = '  -10.32%  '

# Row 47
$dStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'88.57"
$ws.Range("D47").Style = $dStyle
$ws.Range("E47").Value = '  -7.30%  '

# Row 48
$dStyle = $ws.Range("D48").Style
$ws.Range("D48").Value = "'1.340.61"
$ws.Range("D48").Style = $dStyle
$ws.Range("E48").Value = '  -3.90%  '

# Row 49
$dStyle = $ws.Range("D49").Style
$ws.Range("D49").Value = "'7.26"
$ws.Range("D49").Style = $dStyle
$ws.Range("E49").Value = '  -7.10%  '

# Row 50
$ws.Range("E50").Value = '  -4.15%  '

# Row 51
$dStyle = $ws.Range("D51").Style
$ws.Range("D51").Value = "'44.32"
$ws.Range("D51").Style = $dStyle
$ws.Range("E51").Value = '  -6.13%  '
